$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new row of data (row 5)
$ws.Range("A5").Value = "example2"
$ws.Range("B5").Value = "No"
$ws.Range("C5").Value = "Example"
$ws.Range("D5").Value = "Example"
$ws.Range("E5").Value = "Example Work (Year) [url]; Other Work (Year)"
$ws.Range("F5").Value = 50
$ws.Range("G5").Value = 0.5
$ws.Range("H5").Value = 0.8
$ws.Range("I5").Value = "Example Work (Year) [url]; Other Work (Year) [url]"

# Update the selection to match the target (single cell H6)
$ws.Range("H6").Select()
